# SME feedback (Stephanie Ridella, 2023-05-26): the "dated ..." clause
# naming the Power of Attorney date should only appear when
# property_agent_date is actually set, so wrap it in a Jinja {% if %}.
$d = $word.ActiveDocument

$old = "for Property dated {{ property_agent_date }}, empowering"
$new = 'for Property{% if property_agent_date != "" %} dated {{ property_agent_date }}{% endif %} empowering'

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the target 'dated {{ property_agent_date }}' text to replace."
}

Write-Output "found=$found"
